$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.00", "0.994") that must
# remain plain text exactly as scraped, not be coerced into Excel numbers.
# Force Text number-format while assigning, then restore the original style so
# no stray style index is left behind.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '56.028.73'
$ws.Range("E2").Value = '  +3.57%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.314.15'
$ws.Range("E3").Value = '  +2.47%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '516.87'
$ws.Range("E5").Value = '  +4.30%  '

# Row 6
Set-TextValue $ws.Range("D6") '132.74'
$ws.Range("E6").Value = '  +3.33%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.994'

# Row 8
Set-TextValue $ws.Range("D8") '0.534'
$ws.Range("E8").Value = '  +1.93%  '

# Row 9
Set-TextValue $ws.Range("D9") '2.337.73'
$ws.Range("E9").Value = '  +3.41%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.103'
$ws.Range("E10").Value = '  +8.33%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.154'
$ws.Range("E11").Value = '  +1.09%  '

# Row 12
Set-TextValue $ws.Range("D12") '5.15'
$ws.Range("E12").Value = '  +8.52%  '

# Row 13
$ws.Range("E13").Value = '  +2.46%  '

# Row 14
Set-TextValue $ws.Range("D14") '24.08'
$ws.Range("E14").Value = '  +5.97%  '

# Row 15
Set-TextValue $ws.Range("D15") '2.726.56'
$ws.Range("E15").Value = '  +2.55%  '

# Row 16
Set-TextValue $ws.Range("D16") '56.135.30'
$ws.Range("E16").Value = '  +3.83%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.0000135'
$ws.Range("E17").Value = '  +4.69%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.323.95'
$ws.Range("E18").Value = '  +2.87%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.53'
$ws.Range("E19").Value = '  +2.84%  '

# Row 20
Set-TextValue $ws.Range("D20") '4.26'
$ws.Range("E20").Value = '  +2.85%  '

# Row 21
Set-TextValue $ws.Range("D21") '321.63'
$ws.Range("E21").Value = '  +6.71%  '

# Row 22
$ws.Range("E22").Value = '  +5.43%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.999'
$ws.Range("E23").Value = '  -0.13%  '

# Row 24
Set-TextValue $ws.Range("D24") '60.62'
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$ws.Range("E25").Value = '  -1.27%  '

# Row 26
$ws.Range("E26").Value = '  +6.07%  '

# Row 27
Set-TextValue $ws.Range("D27") '7.66'
$ws.Range("E27").Value = '  +5.25%  '

# Row 28
Set-TextValue $ws.Range("D28") '172.42'
$ws.Range("E28").Value = '  +0.99%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D29") '1.68'
$ws.Range("E29").Value = '  +4.90%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D30") '1.18'
$ws.Range("E30").Value = '  +9.08%  '

# Row 31
Set-TextValue $ws.Range("D31") '6.24'
$ws.Range("E31").Value = '  +5.44%  '

# Row 32
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D32") '0.0₃0722'
$ws.Range("E32").Value = '  +4.73%  '

# Row 33
Set-TextValue $ws.Range("D33") '18.36'
$ws.Range("E33").Value = '  +3.48%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.998'

# Row 35
Set-TextValue $ws.Range("D35") '0.991'
$ws.Range("E35").Value = '  -0.57%  '

# Row 36
$ws.Range("E36").Value = '  +5.80%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.927'
$ws.Range("E37").Value = '  -1.84%  '

# Row 38
Set-TextValue $ws.Range("D38") '3.97'
$ws.Range("E38").Value = '  +7.24%  '

# Row 39
$ws.Range("E39").Value = '  +8.85%  '

# Row 40
Set-TextValue $ws.Range("D40") '37.36'
$ws.Range("E40").Value = '  +4.19%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.384'
$ws.Range("E41").Value = '  +3.26%  '

# Row 42
Set-TextValue $ws.Range("D42") '3.66'
$ws.Range("E42").Value = '  +8.98%  '

# Row 43
Set-TextValue $ws.Range("D43") '137.48'
$ws.Range("E43").Value = '  +9.82%  '

# Row 44
Set-TextValue $ws.Range("D44") '5.15'
$ws.Range("E44").Value = '  +7.39%  '

# Row 45
Set-TextValue $ws.Range("D45") '266.87'
$ws.Range("E45").Value = '  +11.18%  '

# Row 46
$ws.Range("E46").Value = '  +3.96%  '

# Row 47
$ws.Range("E47").Value = '  +3.96%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.556'
$ws.Range("E48").Value = '  +2.34%  '

# Row 49
$ws.Range("E49").Value = '  +3.20%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.0215'
$ws.Range("E50").Value = '  +5.77%  '

# Row 51
$ws.Range("E51").Value = '  +5.20%  '
